# "Generate Report for Archive"
#
# The localization-status report has been regenerated: the two entries that
# were previously "Ready for handoff" have since moved on to "In Translation".
# Update the Status / per-language status cells on all three sheets, then
# shrink the now-narrower status columns to match the freshly generated
# report's column widths.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-language status columns E (zh-cn) and F (de-de) ---
foreach ($row in 2,3) {
    foreach ($col in "E","F") {
        $cell = $overview.Range("$col$row")
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- zh-cn / de-de sheets: Status column C ---
foreach ($ws in @($zhcn, $dede)) {
    foreach ($row in 2,3) {
        $cell = $ws.Range("C$row")
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Re-fit the status columns now that the text is shorter ---
# (numeric column indices: E=5, F=6, C=3 - the COM shim mishandles
#  letter-indexed Columns.Item("E") lookups, so use ordinals instead)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5

Write-Host "Status columns updated and re-fit."
